# Scheduled market-data refresh: updates currentAveragePrice/NQ/HQ,
# LevePriceNQ/HQ and LeveProfitNQ/HQ figures across the per-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 63
$ws.Range("H63").Value = 114189.6
$ws.Range("J63").Value = 114189.6
$ws.Range("L63").Value = 114189.6
$ws.Range("N63").Value = -115437.6

# row 66
$ws.Range("H66").Value = 114189.6
$ws.Range("J66").Value = 114189.6
$ws.Range("L66").Value = 342568.8
$ws.Range("N66").Value = -348808.8

# row 112
$ws.Range("H112").Value = 3798.8235
$ws.Range("I112").Value = 1990.5
$ws.Range("K112").Value = 5971.5
$ws.Range("M112").Value = -4863.5

$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 3587.3125
$ws.Range("I45").Value = 2946
$ws.Range("K45").Value = 2946
$ws.Range("M45").Value = -2569

# row 92
$ws.Range("H92").Value = 55708
$ws.Range("J92").Value = 55708
$ws.Range("L92").Value = 55708
$ws.Range("N92").Value = -60700

# row 98
$ws.Range("H98").Value = 21749.75
$ws.Range("J98").Value = 21749.75
$ws.Range("L98").Value = 21749.75
$ws.Range("N98").Value = -27739.75

# row 128
$ws.Range("H128").Value = 90495
$ws.Range("J128").Value = 90495
$ws.Range("L128").Value = 90495
$ws.Range("N128").Value = -100455

$ws = $wb.Worksheets.Item("BSM")
# row 60
$ws.Range("H60").Value = 106943.336
$ws.Range("J60").Value = 106943.336
$ws.Range("L60").Value = 106943.336
$ws.Range("N60").Value = -108141.336

# row 122
$ws.Range("H122").Value = 66663.336
$ws.Range("J122").Value = 49995
$ws.Range("L122").Value = 49995
$ws.Range("N122").Value = -59795

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 3704.3872
$ws.Range("I31").Value = 1995.5625
$ws.Range("K31").Value = 1995.5625
$ws.Range("M31").Value = -1700.5625

# row 34
$ws.Range("H34").Value = 3704.3872
$ws.Range("I34").Value = 1995.5625
$ws.Range("K34").Value = 1995.5625
$ws.Range("M34").Value = -1793.5625

# row 75
$ws.Range("H75").Value = 104637.336
$ws.Range("J75").Value = 104637.336
$ws.Range("L75").Value = 104637.336
$ws.Range("N75").Value = -106633.336

# row 78
$ws.Range("H78").Value = 104637.336
$ws.Range("J78").Value = 104637.336
$ws.Range("L78").Value = 313912.008
$ws.Range("N78").Value = -323896.008

# row 99
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

# row 100
$ws.Range("H100").Value = 103995
$ws.Range("J100").Value = 103995
$ws.Range("L100").Value = 103995
$ws.Range("N100").Value = -106159

# row 126
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

# row 141
$ws.Range("H141").Value = 1173266.8
$ws.Range("J141").Value = 1448292.9
$ws.Range("L141").Value = 1448292.9
$ws.Range("N141").Value = -1458652.9

$ws = $wb.Worksheets.Item("CUL")
# row 9
$ws.Range("H9").Value = 1899.6666
$ws.Range("I9").Value = 349
$ws.Range("K9").Value = 1047
$ws.Range("M9").Value = -823

# row 15
$ws.Range("H15").Value = 48.923077
$ws.Range("I15").Value = 40
$ws.Range("J15").Value = 54.5
$ws.Range("K15").Value = 120
$ws.Range("L15").Value = 163.5
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -443.5

$ws = $wb.Worksheets.Item("GSM")
# row 33
$ws.Range("H33").Value = 23263.334
$ws.Range("I33").Value = 19990
$ws.Range("J33").Value = 24900
$ws.Range("K33").Value = 19990
$ws.Range("L33").Value = 24900
$ws.Range("N33").Value = -25404
$ws.Range("M33").Value = -19738

# row 113
$ws.Range("H113").Value = 2348.8
$ws.Range("I113").Value = 2432.75
$ws.Range("J113").Value = 2013
$ws.Range("K113").Value = 2432.75
$ws.Range("L113").Value = 2013
$ws.Range("M113").Value = -262.75
$ws.Range("N113").Value = -6353

# row 128
$ws.Range("H128").Value = 134679
$ws.Range("J128").Value = 134679
$ws.Range("L128").Value = 134679
$ws.Range("N128").Value = -144639

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 2410
$ws.Range("I22").Value = 1850
$ws.Range("K22").Value = 1850
$ws.Range("M22").Value = -1555

# row 27
$ws.Range("H27").Value = 2410
$ws.Range("I27").Value = 1850
$ws.Range("K27").Value = 1850
$ws.Range("M27").Value = -1743

# row 46
$ws.Range("H46").Value = 2551.75
$ws.Range("I46").Value = 625
$ws.Range("J46").Value = 3194
$ws.Range("K46").Value = 625
$ws.Range("L46").Value = 3194
$ws.Range("M46").Value = -437
$ws.Range("N46").Value = -3570

# row 82
$ws.Range("H82").Value = 12511.6
$ws.Range("I82").Value = 2554.1428
$ws.Range("J82").Value = 17873.309
$ws.Range("K82").Value = 2554.1428
$ws.Range("L82").Value = 17873.309
$ws.Range("M82").Value = -2193.1428
$ws.Range("N82").Value = -18595.309

# row 85
$ws.Range("H85").Value = 12511.6
$ws.Range("I85").Value = 2554.1428
$ws.Range("J85").Value = 17873.309
$ws.Range("K85").Value = 2554.1428
$ws.Range("L85").Value = 17873.309
$ws.Range("M85").Value = -1306.1428
$ws.Range("N85").Value = -20369.309

# row 132
$ws.Range("H132").Value = 4784.6924
$ws.Range("J132").Value = 5666.6665
$ws.Range("L132").Value = 16999.9995
$ws.Range("N132").Value = -22059.9995

$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 4656.8335
$ws.Range("I62").Value = 3630
$ws.Range("J62").Value = 4999.1113
$ws.Range("K62").Value = 3630
$ws.Range("L62").Value = 4999.1113
$ws.Range("M62").Value = -3006
$ws.Range("N62").Value = -6247.1113

# row 64
$ws.Range("H64").Value = 99971.164
$ws.Range("I64").Value = 99950
$ws.Range("K64").Value = 99950
$ws.Range("M64").Value = -99702

# row 65
$ws.Range("H65").Value = 4656.8335
$ws.Range("I65").Value = 3630
$ws.Range("J65").Value = 4999.1113
$ws.Range("K65").Value = 18150
$ws.Range("L65").Value = 24995.5565
$ws.Range("M65").Value = -15030
$ws.Range("N65").Value = -31235.5565

# row 67
$ws.Range("H67").Value = 99971.164
$ws.Range("I67").Value = 99950
$ws.Range("K67").Value = 99950
$ws.Range("M67").Value = -99092

# row 81
$ws.Range("H81").Value = 2521.5
$ws.Range("J81").Value = 3963.1667
$ws.Range("L81").Value = 7926.3334
$ws.Range("N81").Value = -10048.3334

# row 84
$ws.Range("H84").Value = 2521.5
$ws.Range("J84").Value = 3963.1667
$ws.Range("L84").Value = 39631.667
$ws.Range("N84").Value = -50239.667

# row 100
$ws.Range("H100").Value = 1888.6666
$ws.Range("I100").Value = 2239.375
$ws.Range("J100").Value = 1187.25
$ws.Range("K100").Value = 4478.75
$ws.Range("L100").Value = 2374.5
$ws.Range("M100").Value = -3937.75
$ws.Range("N100").Value = -3456.5

# row 109
$ws.Range("H109").Value = 29888
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()

# row 113
$ws.Range("H113").Value = 501.14285
$ws.Range("I113").Value = 676.5
$ws.Range("K113").Value = 2029.5
$ws.Range("M113").Value = 140.5

# row 122
$ws.Range("H122").Value = 3599.3845
$ws.Range("I122").Value = 2900.5557
$ws.Range("J122").Value = 5171.75
$ws.Range("K122").Value = 8701.667099999999
$ws.Range("L122").Value = 15515.25
$ws.Range("M122").Value = -6251.667099999999
$ws.Range("N122").Value = -20415.25

# row 126
$ws.Range("H126").Value = 1043.5625
$ws.Range("I126").Value = 1043.5625
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3130.6875
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -660.6875
$ws.Range("N126").ClearContents()
